# إضافة حدث جديد في Card20: يملأ خلايا الصف 14 الفارغة بـ "nan" (لتتوافق مع
# بقية الصفوف في الورقة) ثم يضيف صف 15 بحدث صيانة جديد.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# Row 14: columns B..K and M were blank placeholder cells; backfill with the
# sheet's literal "nan" placeholder text (matches every other blank cell in
# the table, e.g. row 2's D..O).
$ws.Range("B14:K14").Value = "nan"
$ws.Range("M14").Value = "nan"

# Row 15: new service-log entry for card 20.
# Column A holds a text value "20" (like every other row's card id), so
# force text interpretation via NumberFormat, then drop the resulting
# number-format style so the cell matches its unstyled neighbours (A2..A14).
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "20"
$ws.Range("A15").ClearFormats()

# Columns B..K and M stay blank on the new row (no event data for those
# fields), but are still present as empty-string cell records, consistent
# with the rest of the sheet's blank cells. A leading apostrophe forces an
# empty *text* value instead of an empty/missing cell; ClearFormats then
# drops the transient quote-prefix style.
$ws.Range("B15:K15").Value = "'"
$ws.Range("M15").Value = "'"
$ws.Range("B15:K15").ClearFormats()
$ws.Range("M15").ClearFormats()

$ws.Range("L15").Value = "12\8\2024"
$ws.Range("N15").Value = "تم عمل setting كامل للمكنه وتضيق المسافات"
$ws.Range("O15").Value = "الخبير"
